# PCAR_QTR_FIN.xlsx quarterly financials update
# Inserts two new quarterly columns (2018-12-31 and 2018-09-30) ahead of the
# existing "Period Ending" column D (2018-06-30), shifting the historical
# columns D:K two places to the right (to F:M), and fills in the new
# columns with the newest reported quarter figures. One historical series
# (row 91, "Capital Expenditures") is also restated across the board.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank columns before column D; existing D:K shift to F:M.
$ws.Range("D:E").Insert()

# 2) Rows that get two brand-new values in the new D/E columns while the
#    remaining (now shifted) columns F:M keep their original numbers.
#    Key = row number, Value = @(D-value, E-value)
$newQuarterData = @{
    7   = @(43465, 43373)
    8   = @(6279700, 5756800)
    9   = @(5145900, 4704200)
    10  = @(1133800, 1052600)
    12  = @(80500, 72900)
    14  = @("NA", "NA")
    15  = @(177600, 178500)
    17  = @(5571800, 5112300)
    18  = @(707900, 644500)
    20  = @(44500, 24300)
    21  = @(1012100, 921400)
    23  = @(752400, 668800)
    24  = @(174300, 123500)
    26  = @(578100, 545300)
    27  = @(578100, 545300)
    29  = @("NA", "NA")
    32  = @(-44500, -24300)
    33  = @(578100, 545300)
    35  = @(578100, 545300)
    38  = @(43465, 43373)
    41  = @(3435900, 2914000)
    42  = @(1020400, 1000300)
    43  = @(12155200, 12187600)
    44  = @(1184700, 1292000)
    45  = @(364700, 376800)
    46  = @(18160900, 17770700)
    48  = @(6122500, 5988700)
    52  = @(1199000, 1344000)
    54  = @(25482400, 25103400)
    57  = @(3545800, 3759500)
    59  = @(700200, 4700)
    60  = @(4246000, 3764200)
    61  = @(9950500, 9586900)
    62  = @(2693000, 2581300)
    66  = @(16889500, 15932400)
    72  = @(9275400, 9718400)
    76  = @(8592900, 9171000)
    80  = @(43465, 43373)
    81  = @(578100, 545300)
    83  = @(259700, 252600)
    89  = @(1065100, 728700)
    94  = @(-654200, -531300)
    96  = @(-97700, -98100)
    100 = @(134800, 214900)
    101 = @(-23800, -5100)
    102 = @(521900, 407200)
}

# 3) Rows that are entirely zero across the quarter columns; the new D/E
#    cells need to be filled with 0 to match the rest of the row.
$zeroRows = @(13,22,25,28,30,31,34,47,49,50,51,53,58,63,64,65,68,69,70,71,73,74,75,77,84,85,86,87,88,92,93,97,98,99)

# 4) Rows that are entirely blank across the quarter columns; the new D/E
#    cells stay empty, they just need matching number formatting.
$blankRows = @(11,16,19,39,40,55,56,67,82,90,95)

# Give the freshly inserted D:E cells the same formatting as column F
# (which holds the data that used to live in column D before the insert)
# for every row that actually carries quarterly data.
$allDataRows = @() + $newQuarterData.Keys + $zeroRows + $blankRows
foreach ($r in $allDataRows) {
    $ws.Range("F$r").Copy()
    $ws.Range("D$r:E$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Write the brand-new quarter values.
foreach ($r in $newQuarterData.Keys) {
    $vals = $newQuarterData[$r]
    $ws.Cells.Item($r, 4).Value = $vals[0]
    $ws.Cells.Item($r, 5).Value = $vals[1]
}

# Write zeros for the all-zero rows.
foreach ($r in $zeroRows) {
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 5).Value = 0
}

# 5) Row 91 ("Capital Expenditures") was restated across the whole
#    historical series, not just shifted - overwrite D:M explicitly.
$capex = @(-123100, -112900, -123200, -98400, -127500, -107600, -95600, -429300, -520600, -525100)
for ($i = 0; $i -lt $capex.Length; $i++) {
    $ws.Cells.Item(91, 4 + $i).Value = $capex[$i]
}
